$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.040319800376892
$ws.Range("B1").Value = 0.8059442043304443
$ws.Range("C1").Value = 4.055136680603027
$ws.Range("D1").Value = 2.855939626693726
$ws.Range("E1").Value = 0.7960677146911621
